$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B2 (toplam_hurda_tesviki_butcesi) from 3,000,000,000 to 2,000,000,000
$ws.Range("B2").Value = 2000000000

# Update the active selection on the sheet to F27
$ws.Range("F27").Select()
